$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N64").ClearContents()
$ws.Range("M64").Value = -3352
$ws.Range("J64").Value = 0
$ws.Range("H64").Value = 3600
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 0
$ws.Range("I64").Value = 3600
$ws.Range("M67").Value = -2742
$ws.Range("I67").Value = 3600
$ws.Range("H67").Value = 3600
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("M98").Value = 901.3333
$ws.Range("H98").Value = 585.8461
$ws.Range("K98").Value = 596.6667
$ws.Range("I98").Value = 596.6667
$ws.Range("I122").Value = 596.6667
$ws.Range("M122").Value = 659.9999
$ws.Range("H122").Value = 585.8461
$ws.Range("K122").Value = 1790.0001
$ws.Range("N129").Value = -24237.1
$ws.Range("H129").Value = 3058.318
$ws.Range("J129").Value = 4745.7
$ws.Range("L129").Value = 14237.1
$ws.Range("N133").ClearContents()
$ws.Range("J133").Value = 0
$ws.Range("H133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("I137").Value = 5004.7544
$ws.Range("M137").Value = -12464.2632
$ws.Range("H137").Value = 8743.055
$ws.Range("K137").Value = 15014.2632

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K2").Value = 4020.7368
$ws.Range("L2").Value = 6572.1
$ws.Range("I2").Value = 4020.7368
$ws.Range("N2").Value = -6798.1
$ws.Range("M2").Value = -3907.7368
$ws.Range("J2").Value = 6572.1
$ws.Range("H2").Value = 4900.517
$ws.Range("J61").Value = 24337.666
$ws.Range("H61").Value = 7241.5625
$ws.Range("K61").Value = 3296.3076
$ws.Range("L61").Value = 24337.666
$ws.Range("I61").Value = 3296.3076
$ws.Range("N61").Value = -24761.666
$ws.Range("M61").Value = -3084.3076
$ws.Range("N62").ClearContents()
$ws.Range("J62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("I63").Value = 2290.6086
$ws.Range("M63").Value = -1604.6086
$ws.Range("H63").Value = 2549.3333
$ws.Range("K63").Value = 2290.6086
$ws.Range("H65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("J65").Value = 0
$ws.Range("H66").Value = 2549.3333
$ws.Range("K66").Value = 11453.043
$ws.Range("I66").Value = 2290.6086
$ws.Range("M66").Value = -8021.043
$ws.Range("H110").Value = 2875.7144
$ws.Range("K110").Value = 2624
$ws.Range("I110").Value = 2624
$ws.Range("M110").Value = -579
$ws.Range("H116").Value = 4900.517
$ws.Range("K116").Value = 4020.7368
$ws.Range("L116").Value = 6572.1
$ws.Range("I116").Value = 4020.7368
$ws.Range("N116").Value = -11160.1
$ws.Range("M116").Value = -1726.7368
$ws.Range("J116").Value = 6572.1
$ws.Range("I122").Value = 7640.5293
$ws.Range("M122").Value = -20471.5879
$ws.Range("H122").Value = 7878.9473
$ws.Range("K122").Value = 22921.5879
$ws.Range("K132").Value = 31721.037
$ws.Range("I132").Value = 10573.679
$ws.Range("M132").Value = -29191.037
$ws.Range("H132").Value = 11676.948
$ws.Range("I136").Value = 3296.3076
$ws.Range("M136").Value = -7338.9228
$ws.Range("N136").Value = -78112.99800000001
$ws.Range("H136").Value = 7241.5625
$ws.Range("J136").Value = 24337.666
$ws.Range("K136").Value = 9888.9228
$ws.Range("L136").Value = 73012.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L3").Value = 6572.1
$ws.Range("I3").Value = 4020.7368
$ws.Range("N3").Value = -6800.1
$ws.Range("M3").Value = -3906.7368
$ws.Range("J3").Value = 6572.1
$ws.Range("H3").Value = 4900.517
$ws.Range("K3").Value = 4020.7368
$ws.Range("I16").Value = 1456
$ws.Range("N16").ClearContents()
$ws.Range("M16").Value = -1286
$ws.Range("J16").Value = 0
$ws.Range("H16").Value = 1456
$ws.Range("K16").Value = 1456
$ws.Range("L16").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("I63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("L66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I107").Value = 609.5
$ws.Range("M107").Value = 1310.5
$ws.Range("H107").Value = 884.3077
$ws.Range("K107").Value = 609.5
$ws.Range("H119").Value = 58000
$ws.Range("J119").Value = 58000
$ws.Range("N119").Value = -67676
$ws.Range("L119").Value = 58000
$ws.Range("N122").Value = -17282.8
$ws.Range("J122").Value = 4127.6
$ws.Range("H122").Value = 3233.4
$ws.Range("L122").Value = 12382.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I63").Value = 1111.75
$ws.Range("M63").Value = -2586.25
$ws.Range("H63").Value = 1111.75
$ws.Range("K63").Value = 3335.25
$ws.Range("N64").Value = -64290
$ws.Range("M64").Value = -3927
$ws.Range("J64").Value = 21250
$ws.Range("H64").Value = 9339.4
$ws.Range("K64").Value = 4197
$ws.Range("L64").Value = 63750
$ws.Range("I64").Value = 1399
$ws.Range("H66").Value = 1111.75
$ws.Range("K66").Value = 10005.75
$ws.Range("I66").Value = 1111.75
$ws.Range("M66").Value = -6261.75
$ws.Range("M67").Value = -3261
$ws.Range("I67").Value = 1399
$ws.Range("H67").Value = 9339.4
$ws.Range("J67").Value = 21250
$ws.Range("K67").Value = 4197
$ws.Range("L67").Value = 63750
$ws.Range("N67").Value = -65622
$ws.Range("M75").Value = -2899
$ws.Range("I75").Value = 1299
$ws.Range("J75").Value = 1431.6666
$ws.Range("H75").Value = 1398.5
$ws.Range("K75").Value = 3897
$ws.Range("L75").Value = 4294.9998
$ws.Range("N75").Value = -6290.9998
$ws.Range("I78").Value = 1299
$ws.Range("N78").Value = -22868.9994
$ws.Range("M78").Value = -6699
$ws.Range("J78").Value = 1431.6666
$ws.Range("H78").Value = 1398.5
$ws.Range("K78").Value = 11691
$ws.Range("L78").Value = 12884.9994
$ws.Range("H108").Value = 2869.5
$ws.Range("K108").Value = 8608.5
$ws.Range("I108").Value = 2869.5
$ws.Range("M108").Value = -5728.5
$ws.Range("K140").Value = 5937
$ws.Range("I140").Value = 1979
$ws.Range("M140").Value = -757
$ws.Range("H140").Value = 2014.1724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K2").Value = 76.30768999999999
$ws.Range("I2").Value = 76.30768999999999
$ws.Range("M2").Value = 36.69231000000001
$ws.Range("H2").Value = 100.47059
$ws.Range("K97").Value = 1011.5455
$ws.Range("L97").Value = 2469.6667
$ws.Range("I97").Value = 1011.5455
$ws.Range("N97").Value = -3461.6667
$ws.Range("M97").Value = -515.5454999999999
$ws.Range("J97").Value = 2469.6667
$ws.Range("H97").Value = 1526.1765
$ws.Range("N107").Value = -4896.7142
$ws.Range("I107").Value = 1024.2
$ws.Range("M107").Value = 895.8
$ws.Range("J107").Value = 1056.7142
$ws.Range("H107").Value = 1037.5883
$ws.Range("K107").Value = 1024.2
$ws.Range("L107").Value = 1056.7142
$ws.Range("I122").Value = 1642.0667
$ws.Range("N122").Value = -26255.9995
$ws.Range("M122").Value = -2476.2001
$ws.Range("J122").Value = 7118.6665
$ws.Range("H122").Value = 2554.8333
$ws.Range("K122").Value = 4926.2001
$ws.Range("L122").Value = 21355.9995
$ws.Range("I126").Value = 2075.8333
$ws.Range("M126").Value = -3757.499899999999
$ws.Range("H126").Value = 2045
$ws.Range("K126").Value = 6227.499899999999
$ws.Range("J132").Value = 3732.6667
$ws.Range("K132").Value = 9434.750100000001
$ws.Range("L132").Value = 11198.0001
$ws.Range("I132").Value = 3144.9167
$ws.Range("M132").Value = -6904.750100000001
$ws.Range("N132").Value = -16258.0001
$ws.Range("H132").Value = 3262.4666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M22").Value = -1604.5
$ws.Range("H22").Value = 3499
$ws.Range("K22").Value = 1899.5
$ws.Range("I22").Value = 1899.5
$ws.Range("I27").Value = 1899.5
$ws.Range("M27").Value = -1792.5
$ws.Range("H27").Value = 3499
$ws.Range("K27").Value = 1899.5
$ws.Range("H40").Value = 4995.4287
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4994.6665
$ws.Range("L40").Value = 5000
$ws.Range("I40").Value = 4994.6665
$ws.Range("N40").Value = -5272
$ws.Range("M40").Value = -4858.6665
$ws.Range("J61").Value = 0
$ws.Range("H61").Value = 1499
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("J113").Value = 0
$ws.Range("H113").Value = 1499
$ws.Range("L113").Value = 0
$ws.Range("I122").Value = 3936.5
$ws.Range("N122").Value = -19367.5
$ws.Range("M122").Value = -9359.5
$ws.Range("J122").Value = 4822.5
$ws.Range("H122").Value = 4379.5
$ws.Range("K122").Value = 11809.5
$ws.Range("L122").Value = 14467.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I122").Value = 5725.684
$ws.Range("M122").Value = -14727.052
$ws.Range("H122").Value = 5768.115
$ws.Range("K122").Value = 17177.052
$ws.Range("L126").Value = 1349.5
$ws.Range("I126").Value = 4495.9653
$ws.Range("N126").Value = -8988.5
$ws.Range("M126").Value = -11017.8959
$ws.Range("H126").Value = 4292.968
$ws.Range("J126").Value = 1349.5
$ws.Range("K126").Value = 13487.8959
